$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4747
$ws1.Range("F6").Value = 563
$ws1.Range("F8").Value = 423
$ws1.Range("F21").Value = 64
$ws1.Range("F27").Value = 4208
$ws1.Range("F31").Value = 1999
$ws1.Range("F33").Value = 1958

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4747
$ws4.Range("F6").Value = 563
$ws4.Range("F9").Value = 423
$ws4.Range("F22").Value = 64
$ws4.Range("F28").Value = 4208
$ws4.Range("F34").Value = 1999
$ws4.Range("F36").Value = 1958
